$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A48").Value = "css maths functions part1"
$ws.Range("B48").Value = "https://drive.google.com/file/d/1whT_3lgoF4_o_aCFNDn-wI4gFtSEKR8r/view?usp=sharing"
$ws.Hyperlinks.Add($ws.Range("B48"), "https://drive.google.com/file/d/1whT_3lgoF4_o_aCFNDn-wI4gFtSEKR8r/view?usp=sharing")

$ws.Range("A49").Value = "css maths functions part2"
$ws.Range("B49").Value = "https://drive.google.com/file/d/1rl6UulsflkMEI7jelzhjaa-5QZik36j0/view?usp=sharing"
$ws.Hyperlinks.Add($ws.Range("B49"), "https://drive.google.com/file/d/1rl6UulsflkMEI7jelzhjaa-5QZik36j0/view?usp=sharing")

$ws.Range("B47").Copy()
$ws.Range("B48:B49").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.ActiveWindow.ScrollRow = 37
$ws.Range("B52").Select()
